# Update height and weight jbrelsf2
# (Metadata sheet: bump Version/Status/Date, update Contact info, and
# insert a new "Jurisdiction" property row.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update simple existing values ---
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Insert a new "Jurisdiction" row after the second "Contact" row (row 11) ---
$ws.Rows.Item(12).Insert()
$ws.Range("A13:B13").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# --- Set Description row (now row 13) value ---
$ws.Range("B13").Value = "CBC W Ordered Manual Differential panel - Blood (57782-5)"

# Rows 14 (Purpose) and 15 (Copyright) keep their existing blank values.
# Row 16 (Immutable / BooleanType[null]) was shifted down automatically by the insert.
